$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update SKU ID text for existing rows (the shared strings they pointed to were replaced)
$ws.Range("B2").Value = "FP-W0045-F1Z-2700"
$ws.Range("B3").Value = "FP-43125-KGA-2700"
$ws.Range("B4").Value = "FP-43120-362-2700"

# Update quantity values for existing rows 3 and 4 (from 50 to 20)
$ws.Range("C3").Value = 20
$ws.Range("C4").Value = 20

# Add two new rows (5 and 6) with new SKU IDs, copying formatting from row 4
$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("A5:D5").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("A6:D6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "FP-54410-THU-2700"
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = "Sheet1"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "FP-F533A-RXK-2700"
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = "Sheet1"

# Update the selection to match the target state
$ws.Range("C3").Select() | Out-Null
